# Apply updated symbol list values scraped from coinranking.com
# (GitHub Actions refresh, 2023-01-03 03:49:36 UTC)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h %) hold text that merely looks numeric
# (trailing zeros, percent signs, thousands separators must be preserved verbatim),
# so each write: format the cell as Text, assign the literal string, then restore
# the Normal style so no formatting/style delta is introduced.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

Set-TextValue "D2" '245.75'
Set-TextValue "E2" '1.15%'
Set-TextValue "D3" '29.32'
Set-TextValue "E3" '-0.88%'
Set-TextValue "D4" '5.163'
Set-TextValue "E4" '0.73%'
Set-TextValue "E5" '2.15%'
Set-TextValue "D6" '6.605'
Set-TextValue "E6" '1.68%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D7" '3.164'
Set-TextValue "E7" '5.06%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D8" '0.8582'
Set-TextValue "E8" '3.72%'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D9" '0.8592'
Set-TextValue "E9" '-0.36%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D10" '0.1365'
Set-TextValue "E10" '2.57%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D11" '0.07027'
Set-TextValue "E11" '1.67%'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D12" '0.03294'
Set-TextValue "E12" '1.65%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.02998'
Set-TextValue "E13" '4.96%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.09360'
Set-TextValue "E14" '-0.31%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001541'
Set-TextValue "E15" '2.13%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D16" '0.006011'
Set-TextValue "E16" '-3.28%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.482'
Set-TextValue "E17" '-1.09%'
Set-TextValue "D18" '2.169'
Set-TextValue "E18" '-2.09%'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue "D19" '0.01022'
Set-TextValue "E19" '1,596.63%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue "D20" '0.3202'
Set-TextValue "E20" '1.69%'
Set-TextValue "E21" '-0.97%'
Set-TextValue "D22" '3.318'
Set-TextValue "E22" '-8.23%'
Set-TextValue "D23" '0.04149'
Set-TextValue "E23" '0.18%'
Set-TextValue "D24" '0.1399'
Set-TextValue "E24" '1.93%'
Set-TextValue "E25" '1.33%'
Set-TextValue "E26" '-6.93%'
Set-TextValue "E27" '2.57%'
Set-TextValue "E28" '3.23%'
Set-TextValue "D40" '0.03731'
Set-TextValue "E40" '0.80%'
Set-TextValue "E41" '2.29%'
Set-TextValue "D42" '0.1068'
Set-TextValue "D43" '0.002199'
Set-TextValue "E43" '-4.79%'
Set-TextValue "E44" '-12.56%'
Set-TextValue "D45" '0.00005282'
Set-TextValue "E45" '3.61%'
Set-TextValue "D46" '0.00000000750'
Set-TextValue "E46" '0.03%'
Set-TextValue "D47" '0.05798'
Set-TextValue "E47" '-42.56%'
Set-TextValue "D48" '0.002224'
Set-TextValue "E49" '0.03%'
Set-TextValue "E50" '0.03%'
